# Apply updated "inv_capacity" results from server run across the
# 2025 / 2030 / 2035 sheets (row 2 values).

$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("A2").Value = 0
$ws2025.Range("E2").Value = 0.3836099774358235
$ws2025.Range("G2").Value = 0.2494892361374887
$ws2025.Range("I2").Value = 0.3390794
$ws2025.Range("L2").Value = 0.6371737626639249
$ws2025.Range("M2").Value = 0.07705441666666667
$ws2025.Range("N2").Value = 12.47916728975415
$ws2025.Range("O2").Value = 3.089286314025494

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Value = 0
$ws2030.Range("B2").Value = 0.04315537743582359
$ws2030.Range("E2").Value = 0.2358430747790291
$ws2030.Range("I2").Value = 0.2185782126639248
$ws2030.Range("L2").Value = 0.3134228373360751
$ws2030.Range("M2").Value = 0.04445675000000002
$ws2030.Range("N2").Value = 5.459514395001401
$ws2030.Range("O2").Value = 2.484404990704424

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 0.06436441412223484
$ws2035.Range("B2").Value = 0.04092089999999998
$ws2035.Range("E2").Value = 0.1918809994972733
$ws2035.Range("I2").Value = 0.4873910873360754
$ws2035.Range("L2").Value = 0
$ws2035.Range("M2").Value = 0.03806850000000002
$ws2035.Range("N2").Value = 8.376182333015379
$ws2035.Range("O2").Value = 5.080822541972329
